# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.928.73"
$ws.Range("E2").Value = "  +0.17%  "
$ws.Range("D3").Value = "1.883.19"
$ws.Range("E3").Value = "  +1.27%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "305.25"
$ws.Range("E5").Value = "  +0.22%  "
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("E7").Value = "  +1.87%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3735"
$ws.Range("E8").Value = "  +2.45%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07186"
$ws.Range("E9").Value = "  +0.20%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.04"
$ws.Range("E10").Value = "  +1.71%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8977"
$ws.Range("E11").Value = "  +0.64%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07643"
$ws.Range("E12").Value = "  +2.27%  "
$ws.Range("D13").Value = "1.887.66"
$ws.Range("E13").Value = "  +0.40%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "93.34"
$ws.Range("E14").Value = "  -0.50%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.227"
$ws.Range("E15").Value = "  +0.01%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.000"
$ws.Range("E16").Value = "  -0.02%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008473"
$ws.Range("E17").Value = "  -0.24%  "
$ws.Range("E18").Value = "  +1.33%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9994"
$ws.Range("E19").Value = "  -0.10%  "
$ws.Range("D20").Value = "26.980.71"
$ws.Range("E20").Value = "  +0.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.033"
$ws.Range("E21").Value = "  +0.21%  "
$ws.Range("D22").Value = "2.134.34"
$ws.Range("E22").Value = "  +0.87%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.54"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.379"
$ws.Range("E24").Value = "  -0.59%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.287"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "146.22"
$ws.Range("E26").Value = "  -0.95%  "
$ws.Range("E27").Value = "  -3.28%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.00"
$ws.Range("E28").Value = "  +0.85%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "113.70"
$ws.Range("E29").Value = "  +0.62%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.889"
$ws.Range("E30").Value = "  +4.71%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.761"
$ws.Range("E31").Value = "  +1.58%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09162"
$ws.Range("E32").Value = "  -0.56%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05024"
$ws.Range("E33").Value = "  -1.91%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.227"
$ws.Range("E34").Value = "  +6.61%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7632"
$ws.Range("E35").Value = "  +2.16%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.988"
$ws.Range("E36").Value = "  +0.53%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.258"
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.575"
$ws.Range("E38").Value = "  +0.11%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5567"
$ws.Range("E39").Value = "  +0.19%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01984"
$ws.Range("E40").Value = "  -0.81%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.071"
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.977"
$ws.Range("E42").Value = "  +5.60%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.595"
$ws.Range("E43").Value = "  +0.95%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "118.34"
$ws.Range("E44").Value = "  +1.09%  "
$ws.Range("E45").Value = "  +1.83%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4799"
$ws.Range("E46").Value = "  +2.73%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.9996"
$ws.Range("E47").Value = "  -0.04%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.15"
$ws.Range("E48").Value = "  +1.07%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.593"
$ws.Range("E49").Value = "  +2.20%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "37.54"
$ws.Range("E50").Value = "  +2.34%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "63.76"
$ws.Range("E51").Value = "  +1.34%  "
